$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "POR"
$ws.Range("C2").Value = 14.70909090909091
$ws.Range("B3").Value = "NJN"
$ws.Range("C3").Value = 12.38333333333333
$ws.Range("B4").Value = "CLE"
$ws.Range("C4").Value = 12.95333333333333
$ws.Range("B5").Value = "DAL"
$ws.Range("C5").Value = 15.12142857142857
$ws.Range("B6").Value = "MIA"
$ws.Range("C6").Value = 12.15714285714286
$ws.Range("B7").Value = "SEA"
$ws.Range("C7").Value = 14.31
$ws.Range("B8").Value = "ATL"
$ws.Range("C8").Value = 13.49230769230769
$ws.Range("B9").Value = "MIL"
$ws.Range("C9").Value = 16.00833333333334
$ws.Range("B10").Value = "LAC"
$ws.Range("C10").Value = 11.88461538461539
$ws.Range("B11").Value = "DET"
$ws.Range("C11").Value = 13.28571428571429
$ws.Range("B12").Value = "SAS"
$ws.Range("C12").Value = 12.35714285714286
$ws.Range("B13").Value = "ORL"
$ws.Range("C13").Value = 12.55384615384615
$ws.Range("B14").Value = "UTA"
$ws.Range("C14").Value = 12.61428571428572
$ws.Range("B15").Value = "HOU"
$ws.Range("C15").Value = 13.63076923076923
$ws.Range("B16").Value = "DEN"
$ws.Range("C16").Value = 12.82142857142857
$ws.Range("B17").Value = "LAL"
$ws.Range("C17").Value = 14.09166666666667
$ws.Range("B18").Value = "GSW"
$ws.Range("C18").Value = 11.55333333333333
$ws.Range("B19").Value = "IND"
$ws.Range("C19").Value = 13.31538461538462
$ws.Range("B20").Value = "CHI"
$ws.Range("C20").Value = 15.21666666666667
$ws.Range("B21").Value = "PHI"
$ws.Range("C21").Value = 11.8
$ws.Range("B22").Value = "CHH"
$ws.Range("C22").Value = 12.1
$ws.Range("B23").Value = "BOS"
$ws.Range("C23").Value = 14.1
$ws.Range("B24").Value = "WSB"
$ws.Range("C24").Value = 12.77857142857143
$ws.Range("B25").Value = "SAC"
$ws.Range("C25").Value = 12
$ws.Range("B26").Value = "PHO"
$ws.Range("C26").Value = 19.54615384615385
$ws.Range("B27").Value = "NYK"
$ws.Range("C27").Value = 12.55714285714286
$ws.Range("B28").Value = "MIN"
$ws.Range("C28").Value = 12.37142857142857
